# Refresh the cryptos price/volume table (Price column D, Volume(1h) column E,
# plus a couple of re-ranked rows where Coin name / Link / Price / Volume all
# moved to a different row) to match the latest scrape, per the GitHub Actions
# commit "Updated cryptos list on Fri Mar 10 05:52:28 UTC 2023".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (e.g. "1.000", "273.90") must be forced
# to Text format first, otherwise Excel would coerce them into real numbers
# (e.g. "1.000" -> 1, "273.90" -> 273.9) and the original text formatting
# (trailing zeros, dot-grouping) would be lost.
$textValueUpdates = @{
    'D4' = '1.000'
    'D5' = '1.000'
    'D6' = '273.90'
    'D7' = '0.3722'
    'D9' = '39.82'
    'D10' = '1.012'
    'D11' = '0.06598'
    'D12' = '1.000'
    'D13' = '5.418'
    'D15' = '6.177'
    'D17' = '0.00001007'
    'D18' = '0.05810'
    'D19' = '74.60'
    'D20' = '1.000'
    'D21' = '5.643'
    'D23' = '11.02'
    'D24' = '2.333'
    'D26' = '2.294'
    'D27' = '138.92'
    'D28' = '16.88'
    'D30' = '109.16'
    'D31' = '3.813'
    'D32' = '5.426'
    'D33' = '0.8880'
    'D34' = '0.07741'
    'D35' = '8.449'
    'D36' = '11.30'
    'D37' = '0.05740'
    'D38' = '4.793'
    'D39' = '0.9994'
    'D40' = '0.1926'
    'D41' = '0.02038'
    'D42' = '1.089'
    'D43' = '1.273'
    'D44' = '0.5323'
    'D45' = '3.535'
    'D46' = '12.25'
    'D47' = '0.5131'
    'D48' = '1.798'
    'D49' = '109.70'
    'D50' = '1.050'
    'D51' = '0.9999'
}

foreach ($addr in $textValueUpdates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in $textValueUpdates.Keys) {
    $ws.Range($addr).Value = $textValueUpdates[$addr]
}

# Remaining cells (coin names, links, price strings that already are not
# numeric-parsable, and the volume/percentage strings) can be set directly.
$plainValueUpdates = @{
    'D2' = '20.025.24'
    'E2' = '  -7.98%  '
    'D3' = '1.420.50'
    'E3' = '  -7.73%  '
    'E4' = '  -0.04%  '
    'E5' = '  -0.04%  '
    'E6' = '  -5.45%  '
    'E7' = '  -3.95%  '
    'E8' = '  -3.77%  '
    'E9' = '  -7.62%  '
    'E10' = '  -4.30%  '
    'E11' = '  -8.41%  '
    'E12' = '  -0.05%  '
    'E13' = '  -3.97%  '
    'E14' = '  -7.78%  '
    'E15' = '  -6.38%  '
    'D16' = '1.419.76'
    'E16' = '  -7.87%  '
    'E17' = '  -9.34%  '
    'E18' = '  -11.84%  '
    'E19' = '  -10.46%  '
    'E20' = '  +0.01%  '
    'E21' = '  -8.05%  '
    'E22' = '  -5.75%  '
    'E23' = '  +0.60%  '
    'E24' = '  -3.11%  '
    'D25' = '20.029.78'
    'E25' = '  -7.98%  '
    'E26' = '  -3.38%  '
    'E27' = '  -5.25%  '
    'E28' = '  -8.09%  '
    'D29' = '1.580.93'
    'E29' = '  -7.76%  '
    'E30' = '  -7.12%  '
    'E31' = '  -21.20%  '
    'E32' = '  -8.13%  '
    'E33' = '  -8.47%  '
    'E34' = '  -5.36%  '
    'E35' = '  -4.88%  '
    'E36' = '  +5.89%  '
    'B37' = 'Hedera'
    'C37' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'E37' = '  -5.26%  '
    'B38' = 'InternetComputer(DFINITY)'
    'C38' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'E38' = '  -6.79%  '
    'E39' = '  -0.09%  '
    'E40' = '  -5.50%  '
    'E41' = '  -7.49%  '
    'E42' = '  -8.34%  '
    'E43' = '  -14.21%  '
    'E44' = '  -7.41%  '
    'B45' = 'PancakeSwap'
    'C45' = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    'E45' = '  -5.49%  '
    'B46' = 'EnergySwap'
    'C46' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E46' = '  -5.44%  '
    'E47' = '  -7.06%  '
    'E48' = '  -3.53%  '
    'E49' = '  -6.66%  '
    'E50' = '  -8.22%  '
    'E51' = '  -0.05%  '
}

foreach ($addr in $plainValueUpdates.Keys) {
    $ws.Range($addr).Value = $plainValueUpdates[$addr]
}
